$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "天奇股份"
$ws.Cells.Item(2, 2).Value = "天奇股份"
$ws.Cells.Item(2, 3).Value = "天奇股份"
$ws.Cells.Item(3, 1).Value = "卧龙电驱"
$ws.Cells.Item(3, 2).Value = "利亚德"
$ws.Cells.Item(3, 3).Value = "博纳影业"
$ws.Cells.Item(4, 1).Value = "三花智控"
$ws.Cells.Item(4, 2).Value = "卧龙电驱"
$ws.Cells.Item(4, 3).Value = "风语筑"
$ws.Cells.Item(5, 1).Value = "百达精工"
$ws.Cells.Item(5, 2).Value = "三花智控"
$ws.Cells.Item(5, 3).Value = "嘉美包装"
$ws.Cells.Item(6, 1).Value = "华胜天成"
$ws.Cells.Item(6, 2).Value = "百达精工"
$ws.Cells.Item(6, 3).Value = "光线传媒"
$ws.Cells.Item(7, 1).Value = "博纳影业"
$ws.Cells.Item(7, 2).Value = "中大力德"
$ws.Cells.Item(7, 3).Value = "利欧股份"
$ws.Cells.Item(8, 1).Value = "嘉美包装"
$ws.Cells.Item(8, 2).Value = "绿的谐波"
$ws.Cells.Item(8, 3).Value = "三花智控"
$ws.Cells.Item(9, 1).Value = "五洲新春"
$ws.Cells.Item(9, 2).Value = "嘉美包装"
$ws.Cells.Item(9, 3).Value = "卧龙电驱"
$ws.Cells.Item(10, 1).Value = "光线传媒"
$ws.Cells.Item(10, 2).Value = "鸣志电器"
$ws.Cells.Item(10, 3).Value = "华胜天成"
$ws.Cells.Item(11, 1).Value = "中大力德"
$ws.Cells.Item(11, 2).Value = "长盛轴承"
$ws.Cells.Item(11, 3).Value = "万向钱潮"
$ws.Cells.Item(12, 1).Value = "利欧股份"
$ws.Cells.Item(12, 2).Value = "万向钱潮"
$ws.Cells.Item(12, 3).Value = "百达精工"
$ws.Cells.Item(13, 1).Value = "万向钱潮"
$ws.Cells.Item(13, 2).Value = "金发科技"
$ws.Cells.Item(13, 3).Value = "五洲新春"
$ws.Cells.Item(14, 1).Value = "绿的谐波"
$ws.Cells.Item(14, 2).Value = "光线传媒"
$ws.Cells.Item(14, 3).Value = "汉缆股份"
$ws.Cells.Item(15, 1).Value = "风语筑"
$ws.Cells.Item(15, 2).Value = "利欧股份"
$ws.Cells.Item(15, 3).Value = "掌阅科技"
$ws.Cells.Item(16, 1).Value = "利亚德"
$ws.Cells.Item(16, 2).Value = "贵州茅台"
$ws.Cells.Item(16, 3).Value = "巨力索具"
$ws.Cells.Item(17, 1).Value = "长盛轴承"
$ws.Cells.Item(17, 2).Value = "东方财富"
$ws.Cells.Item(17, 3).Value = "中大力德"
$ws.Cells.Item(18, 1).Value = "克来机电"
$ws.Cells.Item(18, 2).Value = "五洲新春"
$ws.Cells.Item(18, 3).Value = "深科技"
$ws.Cells.Item(19, 1).Value = "金发科技"
$ws.Cells.Item(19, 2).Value = "科大讯飞"
$ws.Cells.Item(19, 3).Value = "绿的谐波"
$ws.Cells.Item(20, 1).Value = "鸣志电器"
$ws.Cells.Item(20, 2).Value = "风语筑"
$ws.Cells.Item(20, 3).Value = "协鑫集成"
$ws.Cells.Item(21, 1).Value = "星环科技-U"
$ws.Cells.Item(21, 2).Value = "华胜天成"
$ws.Cells.Item(21, 3).Value = "长盛轴承"
